# "adjust property of scene"
#
# Scene.xlsx / Sheet1 holds one row per scene (ID, FilePath, MaxGroup, ...,
# CamOffestPos (J), CamOffestRot (K), ...). This change tweaks the camera
# offset position/rotation used for two scenes:
#   - row 2  ("villageScene"): CamOffestPos/CamOffestRot
#   - row 6  ("City")        : CamOffestPos/CamOffestRot
# and leaves the workbook's view focused on the cell that was last edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# villageScene (row 2): CamOffestPos "0,4.2,5.5" -> "0,8,7"
#                        CamOffestRot "25,180"    -> "45,180"
$ws.Range("J2").Value = "0,8,7"
$ws.Range("K2").Value = "45,180"

# City (row 6): CamOffestPos "0,4.2,-5.5" -> "0,8,-7"
#               CamOffestRot "25,0"       -> "45,0"
$ws.Range("J6").Value = "0,8,-7"
$ws.Range("K6").Value = "45,0"

# Scroll column E into view and leave the selection on K7, matching where
# the edit left off.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K7").Select()
